$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114, pushing existing rows 114..124 down to 115..125
$ws.Rows.Item(114).Insert()

# Fill the newly inserted row 114 with the new record's data
$ws.Cells.Item(114, 1).Value = 8
$ws.Cells.Item(114, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(114, 3).Value = "Coquimbo"
$ws.Cells.Item(114, 4).Value = 44714
$ws.Cells.Item(114, 5).Value = 4
$ws.Cells.Item(114, 6).Value = "Fruta"
$ws.Cells.Item(114, 7).Value = 100109
$ws.Cells.Item(114, 8).Value = "Uva"
$ws.Cells.Item(114, 9).Value = 100109001
$ws.Cells.Item(114, 10).Value = "Uva"
$ws.Cells.Item(114, 11).Value = "Red Globe"
$ws.Cells.Item(114, 12).Value = "Primera"
$ws.Cells.Item(114, 13).Value = 300
$ws.Cells.Item(114, 14).Value = 8000
$ws.Cells.Item(114, 15).Value = 9000
$ws.Cells.Item(114, 16).Value = 8500
$ws.Cells.Item(114, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(114, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(114, 19).Value = 472
$ws.Cells.Item(114, 20).Value = 18

# Match the date formatting style used by column D in the rest of the table
$ws.Cells.Item(114, 4).NumberFormat = $ws.Cells.Item(115, 4).NumberFormat
